$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown Chart")
$ws.Range("B20").Value = 4
$ws.Range("B21").Value = 3
$ws.Range("B22").Value = 1
